$d = $word.ActiveDocument

$d.Content.Find.Execute("Want to invoke fast paced speed running feeling.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Want to invoke fast paced speed running feeling.^pThis will be a Windowed Desktop Game.", 2)
